$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.431.39'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '2.532.62'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.82%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.65'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.68%  '
$ws.Range('D16').Value = '2.533.26'
$ws.Range('E16').Value = '  -2.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.830'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').Value = '42.435.47'
$ws.Range('E18').Value = '  -1.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.42%  '
$ws.Range('D20').Value = '0.0₃0948'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.95'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '243.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.95%  '
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.30'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.36'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.70'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +15.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0795'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E35').Value = '  -3.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.02'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.53%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.16'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.64%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.14'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.82%  '
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.60%  '
$ws.Range('E42').Value = '  -3.78%  '
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0296'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('D46').Value = '1.965.76'
$ws.Range('E46').Value = '  -1.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('D48').Value = '2.773.68'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '80.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.192'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '101.50'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.67%  '
